$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

$ws.Range("B2").Value = 90.83900850744227
$ws.Range("C2").Value = 90.15903210483022
$ws.Range("D2").Value = 88.31700381052991
$ws.Range("E2").Value = 91.1467091320265

$ws.Range("B3").Value = 98.37587613993378
$ws.Range("C3").Value = 97.30243364099044
$ws.Range("D3").Value = 98.22366108225965
$ws.Range("E3").Value = 98.00312514745758

$ws.Range("B4").Value = 99.25189821715138
$ws.Range("C4").Value = 99.17676852637356
$ws.Range("D4").Value = 99.27027355081542
$ws.Range("E4").Value = 99.28358666821494

$ws.Range("B5").Value = 98.73279400035089
$ws.Range("C5").Value = 98.74942663539709
$ws.Range("D5").Value = 98.73160709071686
$ws.Range("E5").Value = 98.71779344194712

$ws.Range("B6").Value = 98.32670240854404
$ws.Range("C6").Value = 98.22815838751862
$ws.Range("D6").Value = 98.25168058376293
$ws.Range("E6").Value = 98.20107411073764

$ws.Range("B7").Value = 97.26390962487595
$ws.Range("C7").Value = 97.26158248607535
$ws.Range("D7").Value = 97.30807456446634
$ws.Range("E7").Value = 97.29216094839659

$ws.Range("B8").Value = 95.90080237464352
$ws.Range("C8").Value = 95.86813463132646
$ws.Range("D8").Value = 95.88320732434333
$ws.Range("E8").Value = 95.83807716395438
